$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to remain plain text so numeric-looking price strings
# (e.g. "1.00", "0.999") are not silently converted to numbers, which would
# drop trailing zeros / change their textual representation.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '65.543.46'
$ws.Range("E2").Value = '  +2.70%  '
$ws.Range("D3").Value = '3.485.38'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '581.65'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").Value = '161.08'
$ws.Range("E6").Value = '  +2.73%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '0.609'
$ws.Range("E8").Value = '  +11.53%  '
$ws.Range("D9").Value = '3.489.90'
$ws.Range("E9").Value = '  +1.68%  '
$ws.Range("D10").Value = '7.30'
$ws.Range("E10").Value = '  -2.00%  '
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +1.95%  '
$ws.Range("E12").Value = '  +1.68%  '
$ws.Range("D13").Value = '4.093.59'
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = '0.0000195'
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("D16").Value = '28.71'
$ws.Range("E16").Value = '  +5.29%  '
$ws.Range("D17").Value = '65.553.64'
$ws.Range("E17").Value = '  +2.48%  '
$ws.Range("D18").Value = '3.486.65'
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("D19").Value = '6.45'
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("D20").Value = '14.30'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = '388.78'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = '8.24'
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").Value = '0.552'
$ws.Range("E23").Value = '  +2.45%  '
$ws.Range("D24").Value = '73.31'
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  +2.81%  '
$ws.Range("D27").Value = '10.15'
$ws.Range("E27").Value = '  +5.99%  '
$ws.Range("E28").Value = '  +0.63%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '6.33'
$ws.Range("E30").Value = '  +4.11%  '
$ws.Range("D31").Value = '1.44'
$ws.Range("E31").Value = '  +5.08%  '
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  +2.49%  '
$ws.Range("D33").Value = '23.68'
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '7.18'
$ws.Range("E34").Value = '  +3.67%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  +6.63%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '162.66'
$ws.Range("E36").Value = '  +2.14%  '
$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").Value = '1.95'
$ws.Range("E37").Value = '  +5.73%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '3.064.09'
$ws.Range("E38").Value = '  +4.83%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.0774'
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '27.20'
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0322'
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '4.56'
$ws.Range("E42").Value = '  +3.06%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '43.13'
$ws.Range("E43").Value = '  +3.89%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '6.53'
$ws.Range("E44").Value = '  +0.78%  '
$ws.Range("D45").Value = '0.779'
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("D46").Value = '25.92'
$ws.Range("E46").Value = '  +9.40%  '
$ws.Range("D47").Value = '1.12'
$ws.Range("E47").Value = '  +3.02%  '
$ws.Range("D48").Value = '317.70'
$ws.Range("E48").Value = '  +8.45%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.110'
$ws.Range("E49").Value = '  +5.88%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '2.22'
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '6.73'
$ws.Range("E51").Value = '  +3.59%  '
